$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new cell text in the order that reproduces the target shared-string table order:
# D8, D9, D7, C9, B8, B9
$ws.Range("D8").Value = "MongoDB locally adapted (80%)."
$ws.Range("D9").Value = "MongoDB installation aphrodite ""finished"" (Mr. Oberlercher adapts aphrodite) (99%). Research and queries for the collections (40%)"
$ws.Range("D7").Value = "MongoDB in Virtual Machine installation finished  (60%)"
$ws.Range("C9").Value = "C# GoogleMaps WebBrowser dynamically added. 2nd column in progress (70%)."
$ws.Range("B8").Value = "Android Research connection MongoDB with AndroidStudio (20%)"
$ws.Range("B9").Value = "Android Creation of remaining activities, MongoDB jar added, connection in progress (40%)"

# Update the active cell selection to B9
$ws.Range("B9").Select()
